# Update the capital-structure / fundamentals database for the Australia
# Insurance (Prop/Cas.) industry group: refresh metric values for rows 2-5
# and swap the company order of rows 4 and 5 (Suncorp / QBE).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -0.0378
$ws.Range("E2").Value = -0.0701
$ws.Range("F2").Value = 0.0687
$ws.Range("G2").Value = 0.1132008688556962
$ws.Range("H2").Value = 0.1132008688556962
$ws.Range("I2").Value = 0.05585341227053751
$ws.Range("J2").Value = 0.04861788949090276
$ws.Range("K2").Value = 305.7
$ws.Range("L2").Value = 0.01138975927630133
$ws.Range("M2").Value = 1608.8
$ws.Range("N2").Value = 0.0573189632136813
$ws.Range("O2").Value = 5.262675825973177
$ws.Range("P2").Value = 1513.8
$ws.Range("Q2").Value = 0.05393426560969092
$ws.Range("R2").Value = 4.951913640824338
$ws.Range("S2").Value = 95
$ws.Range("T2").Value = 0.05905022376926902
$ws.Range("U2").Value = 3536.8
$ws.Range("V2").Value = 0.1260105103767703
$ws.Range("W2").Value = 0.0668774914816382
$ws.Range("X2").Value = 0.0476326950273929
$ws.Range("Y2").Value = 0.0192447964542453
$ws.Range("Z2").Value = 0.7413231175460014
$ws.Range("AA2").Value = 0.0432295483060257
$ws.Range("AB2").Value = 0.03945195722570763
$ws.Range("AC2").Value = 0.005427343788318884
$ws.Range("AD2").Value = 21225.4
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 21225.4
$ws.Range("AG2").Value = 17688.6
$ws.Range("AH2").Value = 0.4305975099862252
$ws.Range("AI2").Value = 0.4924184074015646
$ws.Range("AJ2").Value = 0.3865845209709743
$ws.Range("AK2").Value = 0.4470475843872259
$ws.Range("AL2").Value = 1006.2
$ws.Range("AM2").Value = 1006.2
$ws.Range("AN2").Value = 12.35255776057731
$ws.Range("AO2").Value = 1.489862850327967
$ws.Range("AP2").Value = 10.29424431123785
$ws.Range("AQ2").Value = 1.489862850327967
$ws.Range("D3").Value = -0.0828
$ws.Range("E3").Value = -0.09789999999999999
$ws.Range("F3").Value = -0.0467
$ws.Range("G3").Value = 0.09171164302481892
$ws.Range("H3").Value = 0.09171164302481892
$ws.Range("I3").Value = 0.0802150839670151
$ws.Range("J3").Value = 0.07476917171685833
$ws.Range("K3").Value = 300.3
$ws.Range("L3").Value = 0.06025160008828073
$ws.Range("M3").Value = 478.5
$ws.Range("N3").Value = 0.05441086170430511
$ws.Range("O3").Value = 1.593406593406593
$ws.Range("P3").Value = 478.5
$ws.Range("Q3").Value = 0.05441086170430511
$ws.Range("R3").Value = 1.593406593406593
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 1586.6
$ws.Range("V3").Value = 0.1804143640126447
$ws.Range("W3").Value = 0.0668774914816382
$ws.Range("X3").Value = 0.04390589336435248
$ws.Range("Y3").Value = 0.02297159811728572
$ws.Range("Z3").Value = 1.794067888124977
$ws.Range("AA3").Value = 0.1341409699989178
$ws.Range("AB3").Value = 0.03945195722570763
$ws.Range("AC3").Value = 0.09468901277321018
$ws.Range("AD3").Value = 1505.8
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 1505.8
$ws.Range("AG3").Value = -80.79999999999995
$ws.Range("AH3").Value = 0.1461941747572816
$ws.Range("AI3").Value = 0.2460256514990605
$ws.Range("AJ3").Value = -0.009273073656666736
$ws.Range("AK3").Value = -0.01782130174904607
$ws.Range("AL3").Value = 63.5
$ws.Range("AM3").Value = 63.5
$ws.Range("AN3").Value = 3.294246335593962
$ws.Range("AO3").Value = 6.296062992125984
$ws.Range("AP3").Value = -0.1767665718661123
$ws.Range("AQ3").Value = 6.296062992125984
$ws.Range("B4").Value = "Suncorp Group Limited (ASX:SUN)"
$ws.Range("D4").Value = -0.0245
$ws.Range("E4").Value = -0.0423
$ws.Range("F4").Value = 0.0832
$ws.Range("G4").Value = 0.1458496783304566
$ws.Range("H4").Value = 0.1458496783304566
$ws.Range("I4").Value = 0.1394849364506512
$ws.Range("J4").Value = 0.09474608842063728
$ws.Range("K4").Value = 630.4
$ws.Range("L4").Value = 0.06182331711909619
$ws.Range("M4").Value = 619.3
$ws.Range("N4").Value = 0.06455143372351181
$ws.Range("O4").Value = 0.9823921319796954
$ws.Range("P4").Value = 619.3
$ws.Range("Q4").Value = 0.06455143372351181
$ws.Range("R4").Value = 0.9823921319796954
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 1012.2
$ws.Range("V4").Value = 0.1055045393427074
$ws.Range("W4").Value = 0.06856347341860262
$ws.Range("X4").Value = 0.07716701052780199
$ws.Range("Y4").Value = -0.008603537109199366
$ws.Range("Z4").Value = 0.4562673670927989
$ws.Range("AA4").Value = 0.0432295483060257
$ws.Range("AB4").Value = 0.03780220451770681
$ws.Range("AC4").Value = 0.005427343788318884
$ws.Range("AD4").Value = 16394.6
$ws.Range("AF4").Value = 16394.6
$ws.Range("AG4").Value = 15382.4
$ws.Range("AH4").Value = 0.6308405640956576
$ws.Range("AI4").Value = 0.6500402444004774
$ws.Range("AJ4").Value = 0.6158798541016884
$ws.Range("AK4").Value = 0.6354079318592076
$ws.Range("AL4").Value = 689.7
$ws.Range("AM4").Value = 689.7
$ws.Range("AN4").Value = 11.02380311995696
$ws.Range("AO4").Value = 2.062200956937799
$ws.Range("AP4").Value = 10.34319526627219
$ws.Range("AQ4").Value = 2.062200956937799
$ws.Range("B5").Value = "QBE Insurance Group Limited (ASX:QBE)"
$ws.Range("D5").Value = -0.0378
$ws.Range("E5").ClearContents()
$ws.Range("F5").Value = 0.0687
$ws.Range("G5").Value = 0.0938330903164937
$ws.Range("H5").Value = 0.0938330903164937
$ws.Range("I5").Value = -0.02770391971867227
$ws.Range("J5").Value = -0.02770391971867227
$ws.Range("K5").Value = -625
$ws.Range("L5").Value = -0.05360665580238443
$ws.Range("M5").Value = 511
$ws.Range("N5").Value = 0.05279252846250801
$ws.Range("O5").Value = -0.8176
$ws.Range("P5").Value = 416
$ws.Range("Q5").Value = 0.04297787052916503
$ws.Range("R5").Value = -0.6656
$ws.Range("S5").Value = 95
$ws.Range("T5").Value = 0.1859099804305284
$ws.Range("U5").Value = 938
$ws.Range("V5").Value = 0.09690683306816539
$ws.Range("W5").Value = -0.07470714797991872
$ws.Range("X5").Value = 0.0476326950273929
$ws.Range("Y5").Value = -0.1223398430073116
$ws.Range("Z5").Value = 1.052351295243253
$ws.Range("AA5").Value = -0.02915425579925986
$ws.Range("AB5").Value = 0.03950471638636427
$ws.Range("AC5").Value = -0.06865897218562413
$ws.Range("AD5").Value = 3325
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 3325
$ws.Range("AG5").Value = 2387
$ws.Range("AH5").Value = 0.2556826920119344
$ws.Range("AI5").Value = 0.2826659865680524
$ws.Range("AJ5").Value = 0.1978220513160512
$ws.Range("AK5").Value = 0.2205080831408776
$ws.Range("AL5").Value = 253
$ws.Range("AM5").Value = 253
$ws.Range("AN5").Value = -14.71238938053097
$ws.Range("AO5").Value = -1.276679841897233
$ws.Range("AP5").Value = -10.56194690265487
$ws.Range("AQ5").Value = -1.276679841897233
